# Applies review-record baseline update to ORM_Records workbook:
#  1. Rename Sheet1 -> "Review Records"
#  2. Move the active selection to I11
#  3. Set the "Closed Date" (F16) for the STP.docx review row to 2019-01-05
#     (baseline date after the review comments were incorporated)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the first (active) worksheet
$ws.Name = "Review Records"

# 2. Update the selected / active cell on that sheet
$ws.Activate()
$ws.Range("I11").Select()

# 3. Set the Closed Date value for row 16 (STP.docx) and match the date
#    formatting already used by the Issue Date column (E16)
[DateTime]$closedDate = "2019-01-05"
$ws.Range("F16").Value = $closedDate
